$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Year/mean_ETR/max_ETR statistics for rows 6-16 (Sheet 1).
# Rows 6-13: years were off by one (shifted), correct to true year;
#            rows 10-13 additionally get corrected mean/max ETR values.
# Rows 14-16: mean/max ETR values corrected (year/month unchanged).

$ws.Range("A6").Value = 2021
$ws.Range("A7").Value = 2021
$ws.Range("C7").Value = 15.580376344086
$ws.Range("D7").Value = 141.7

$ws.Range("A8").Value = 2022
$ws.Range("C8").Value = 18.5020161290323
$ws.Range("D8").Value = 129.1

$ws.Range("A9").Value = 2022
$ws.Range("C9").Value = 2.48020833333335
$ws.Range("D9").Value = 48.6

$ws.Range("A10").Value = 2020
$ws.Range("C10").Value = 0.00381482281111117
$ws.Range("D10").Value = 0.6

$ws.Range("A11").Value = 2020
$ws.Range("C11").Value = 47.4453975682084
$ws.Range("D11").Value = 361.566

$ws.Range("A12").Value = 2021
$ws.Range("C12").Value = 41.2198884408602
$ws.Range("D12").Value = 373.394

$ws.Range("A13").Value = 2021
$ws.Range("C13").Value = 0.5425431547619
$ws.Range("D13").Value = 46.54

$ws.Range("C14").Value = 0.00879305555555557
$ws.Range("D14").Value = 1.49

$ws.Range("C15").Value = 26.4007096774193
$ws.Range("D15").Value = 489.998

$ws.Range("C16").Value = 12.6349986559139
$ws.Range("D16").Value = 465.6
